# Scheduled market-data refresh for the Leve profit tables.
# Updates cached currentAveragePrice / LevePrice / LeveProfit columns (H:N)
# per sheet+row based on the latest pull; some rows gain or lose their
# LeveProfitNQ (M) / LeveProfitHQ (N) cell depending on whether a profit
# value is present for that Leve.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 15873855
$ws.Range("I19").Value = 35714960
$ws.Range("J19").Value = 971.6
$ws.Range("K19").Value = 35714960
$ws.Range("L19").Value = 971.6
$ws.Range("M19").Value = -35714785
$ws.Range("N19").Value = -1321.6

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("N44").ClearContents()

$ws.Range("H74").Value = 4365.7144
$ws.Range("I74").Value = 4000
$ws.Range("J74").Value = 4426.6665
$ws.Range("K74").Value = 4000
$ws.Range("L74").Value = 4426.6665
$ws.Range("M74").Value = -3064
$ws.Range("N74").Value = -6298.6665

$ws.Range("H75").Value = 28525.8
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 28525.8
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 28525.8
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -30397.8

$ws.Range("H77").Value = 4365.7144
$ws.Range("I77").Value = 4000
$ws.Range("J77").Value = 4426.6665
$ws.Range("K77").Value = 20000
$ws.Range("L77").Value = 22133.3325
$ws.Range("M77").Value = -15320
$ws.Range("N77").Value = -31493.3325

$ws.Range("H78").Value = 28525.8
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 28525.8
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 85577.39999999999
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -94937.39999999999

$ws.Range("H141").Value = 563217.4399999999
$ws.Range("I141").Value = 1726.6666
$ws.Range("J141").Value = 1766412
$ws.Range("K141").Value = 5179.9998
$ws.Range("L141").Value = 5299236
$ws.Range("M141").Value = 0.0002000000004045432
$ws.Range("N141").Value = -5309596

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2221.303
$ws.Range("I61").Value = 1270.4814
$ws.Range("J61").Value = 6500
$ws.Range("K61").Value = 1270.4814
$ws.Range("L61").Value = 6500
$ws.Range("M61").Value = -1058.4814
$ws.Range("N61").Value = -6924

$ws.Range("H132").Value = 25003432
$ws.Range("I132").Value = 58825824
$ws.Range("J132").Value = 4273.5654
$ws.Range("K132").Value = 176477472
$ws.Range("L132").Value = 12820.6962
$ws.Range("M132").Value = -176474942
$ws.Range("N132").Value = -17880.6962

$ws.Range("H136").Value = 2221.303
$ws.Range("I136").Value = 1270.4814
$ws.Range("J136").Value = 6500
$ws.Range("K136").Value = 3811.4442
$ws.Range("L136").Value = 19500
$ws.Range("M136").Value = -1261.4442
$ws.Range("N136").Value = -24600

$ws.Range("H137").Value = 29642.857
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 29642.857
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 29642.857
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -39842.857

$ws.Range("H139").Value = 29547.37
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 29547.37
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 29547.37
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -39827.37

$ws.Range("H141").Value = 29498.334
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 29498.334
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 29498.334
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -39858.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 30000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 30000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 30000
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -40120

$ws.Range("H134").Value = 2303.4167
$ws.Range("I134").Value = 1434.421
$ws.Range("J134").Value = 5605.6
$ws.Range("K134").Value = 4303.263
$ws.Range("L134").Value = 16816.8
$ws.Range("M134").Value = -1768.263
$ws.Range("N134").Value = -21886.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 44001.2
$ws.Range("I4").Value = 5000
$ws.Range("J4").Value = 70002
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 70002
$ws.Range("M4").Value = -4888
$ws.Range("N4").Value = -70226

$ws.Range("H99").Value = 10000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 10000
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -12996

$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 30000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -34940

$ws.Range("H132").Value = 2424.2778
$ws.Range("I132").Value = 1770.5
$ws.Range("J132").Value = 4712.5
$ws.Range("K132").Value = 5311.5
$ws.Range("L132").Value = 14137.5
$ws.Range("M132").Value = -2781.5
$ws.Range("N132").Value = -19197.5

$ws.Range("H134").Value = 1653.2
$ws.Range("I134").Value = 896.35
$ws.Range("J134").Value = 2662.3333
$ws.Range("K134").Value = 2689.05
$ws.Range("L134").Value = 7986.999899999999
$ws.Range("M134").Value = -154.0500000000002
$ws.Range("N134").Value = -13056.9999

$ws.Range("H135").Value = 28082.857
$ws.Range("I135").Value = 20000
$ws.Range("J135").Value = 29430
$ws.Range("K135").Value = 20000
$ws.Range("L135").Value = 29430
$ws.Range("M135").Value = -14930
$ws.Range("N135").Value = -39570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 2337.4
$ws.Range("I47").Value = 391.5
$ws.Range("J47").Value = 3634.6667
$ws.Range("K47").Value = 1174.5
$ws.Range("L47").Value = 10904.0001
$ws.Range("M47").Value = -743.5
$ws.Range("N47").Value = -11766.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10285.571
$ws.Range("I5").Value = 6333.3335
$ws.Range("J5").Value = 13249.75
$ws.Range("K5").Value = 6333.3335
$ws.Range("L5").Value = 13249.75
$ws.Range("M5").Value = -6221.3335
$ws.Range("N5").Value = -13473.75

$ws.Range("H88").Value = 29681.818
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 29681.818
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 29681.818
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -30583.818

$ws.Range("H91").Value = 29681.818
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 29681.818
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 29681.818
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -32801.818

$ws.Range("H132").Value = 3734
$ws.Range("I132").Value = 3534.6667
$ws.Range("J132").Value = 3973.2
$ws.Range("K132").Value = 10604.0001
$ws.Range("L132").Value = 11919.6
$ws.Range("M132").Value = -8074.000100000001
$ws.Range("N132").Value = -16979.6

$ws.Range("H137").Value = 29642.857
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 29642.857
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 29642.857
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -39842.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("N124").ClearContents()

$ws.Range("H136").Value = 3752.476
$ws.Range("I136").Value = 3446.8
$ws.Range("J136").Value = 4516.6665
$ws.Range("K136").Value = 10340.4
$ws.Range("L136").Value = 13549.9995
$ws.Range("M136").Value = -7790.400000000001
$ws.Range("N136").Value = -18649.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 4845
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 4845
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 4845
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -5425

$ws.Range("H80").Value = 40000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 40000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 40000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -41996

$ws.Range("H83").Value = 40000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 40000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 120000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -129984

$ws.Range("H122").Value = 1670834
$ws.Range("I122").Value = 3335334.8
$ws.Range("J122").Value = 6333.3335
$ws.Range("K122").Value = 10006004.4
$ws.Range("L122").Value = 19000.0005
$ws.Range("M122").Value = -10003554.4
$ws.Range("N122").Value = -23900.0005

$ws.Range("H132").Value = 13627.659
$ws.Range("I132").Value = 2113.8333
$ws.Range("J132").Value = 21598.77
$ws.Range("K132").Value = 6341.499899999999
$ws.Range("L132").Value = 64796.31
$ws.Range("M132").Value = -3811.499899999999
$ws.Range("N132").Value = -69856.31

$ws.Range("H136").Value = 1395.1154
$ws.Range("I136").Value = 864.5333000000001
$ws.Range("J136").Value = 2118.6365
$ws.Range("K136").Value = 2593.5999
$ws.Range("L136").Value = 6355.9095
$ws.Range("M136").Value = -43.59990000000016
$ws.Range("N136").Value = -11455.9095
